$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (symbol list update).
# Source values are stored as literal text (e.g. "302.25", "1.70%"),
# not as numbers/percentages, so force Text format before writing
# to avoid Excel re-interpreting numeric-looking strings.
$cells = @{
    'D2' = '302.25'
    'E2' = '1.70%'
    'D3' = '32.11'
    'E3' = '2.58%'
    'D4' = '5.131'
    'E4' = '0.80%'
    'D5' = '0.07913'
    'E5' = '-0.97%'
    'D6' = '2.284'
    'E6' = '-7.07%'
    'D7' = '7.829'
    'E7' = '0.47%'
    'D8' = '3.806'
    'E8' = '0.09%'
    'D9' = '0.9291'
    'E9' = '1.00%'
    'D10' = '0.1776'
    'E10' = '2.81%'
    'D11' = '0.07726'
    'E11' = '6.12%'
    'D12' = '0.08876'
    'E12' = '1.57%'
    'D13' = '0.03091'
    'E13' = '1.94%'
    'E14' = '0.39%'
    'D15' = '0.001506'
    'E15' = '0.79%'
    'D16' = '0.005965'
    'E16' = '0.58%'
    'D17' = '3.475'
    'E17' = '-0.74%'
    'D18' = '2.254'
    'E18' = '0.35%'
    'E19' = '0.24%'
    'D20' = '0.1342'
    'E20' = '0.31%'
    'D21' = '4.241'
    'E21' = '-7.85%'
    'E22' = '10.68%'
    'D23' = '0.04592'
    'E23' = '-1.19%'
    'D24' = '0.001248'
    'E24' = '0.06%'
    'D25' = '0.004512'
    'E25' = '1.83%'
    'D26' = '0.0001246'
    'E26' = '3.65%'
    'D39' = '0.01796'
    'E39' = '0.35%'
    'D40' = '0.04799'
    'E40' = '7.87%'
    'D41' = '0.007294'
    'E41' = '5.12%'
    'D42' = '0.1374'
    'E42' = '2.38%'
    'D43' = '0.002117'
    'E43' = '-2.18%'
    'D44' = '0.01095'
    'E44' = '11.50%'
    'D45' = '0.00006237'
    'E45' = '-5.52%'
    'D46' = '0.00000000748'
    'E46' = '-0.39%'
    'D47' = '0.002502'
    'E47' = '-52.23%'
    'D48' = '0.7066'
    'E48' = '-13.89%'
    'D49' = '0.00002094'
    'E49' = '-0.39%'
    'D50' = '0.0001994'
    'E50' = '-0.39%'
}

foreach ($ref in $cells.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$ref]
}
